{"js": "// Locate the \"Potential stakeholders ...\" paragraph and turn it into the\n// new \"Myself: / Engineer: / Apprentice:\" stakeholder list.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"Potential stakeholders\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) === 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Potential stakeholders' paragraph\");\n}\n\n// Replace the whole paragraph's text with \"Myself:\" then append the\n// explanatory sentence, producing:\n//   Myself: Since I am in charge of the development of system, any hold\n//   ups or push backs directly affect me.\ntarget.insertText(\"Myself:\", \"Replace\");\ntarget.insertText(\n  \" Since I am in charge of the development of system, any hold ups or push backs directly affect me. \",\n  \"End\"\n);\n\n// Add the two follow-up paragraphs right after it.\nconst engineerPara = target.insertParagraph(\"Engineer:\", \"After\");\nengineerPara.insertParagraph(\"Apprentice:\", \"After\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Find the paragraph that starts the \"Identification of stakeholders\"\n# answer (\"Potential stakeholders for this project include myself, ...\")\n$targetIndex = 0\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text.StartsWith(\"Potential stakeholders\")) {\n        $targetIndex = $i\n        break\n    }\n}\n\n$target = $d.Paragraphs.Item($targetIndex)\n$r = $target.Range\n$null = $r.Find.Execute(\n    \"Potential stakeholders for this project include myself, engineers designing structures with small areas and \",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Myself:\", 2\n)\n\n# $r now spans the freshly-inserted \"Myself:\" \u2014 append the explanation right\n# after it, inside the same paragraph.\n$r.InsertAfter(\" Since I am in charge of the development of system, any hold ups or push backs directly affect me. \")\n\n# Add the two follow-up paragraphs right after this one.\n$p1 = $d.Paragraphs.Item($targetIndex)\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Item($targetIndex + 1)\n$p2.Range.Text = \"Engineer:\"\n$p2.Range.InsertParagraphAfter()\n$p3 = $d.Paragraphs.Item($targetIndex + 2)\n$p3.Range.Text = \"Apprentice:\"\n"}
